$wb = $excel.ActiveWorkbook

# Locate the "Czech" sheet - it is the template for the new "Swiss" market sheet.
$czech = $wb.Worksheets.Item("Czech")

# Duplicate the Czech sheet to the end of the workbook and rename it "Swiss".
$czech.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Swiss"

# Update the market description and user-story reference for Switzerland.
$newSheet.Range("B2").Value = "Switzerland Market"
$newSheet.Range("B4").Value = "NGC-3476/T2653"

# Rename the existing P32AR / P32DR entries to their Swiss (CH) variants.
$newSheet.Range("A16").Value = "P32AR-CH"
$newSheet.Range("A17").Value = "P32DR-CH"

# Insert a new row for PR1DSCH right after the PR1DS row, copying its formatting.
$newSheet.Rows.Item(19).Insert()
$newSheet.Range("A18").Copy()
$newSheet.Range("A19").PasteSpecial(-4122)
$newSheet.Range("A19").Value = "PR1DSCH"

# Insert a new row for PR8ASCH right after the PR8AS row, copying its formatting.
$newSheet.Rows.Item(21).Insert()
$newSheet.Range("A20").Copy()
$newSheet.Range("A21").PasteSpecial(-4122)
$newSheet.Range("A21").Value = "PR8ASCH"

$excel.CutCopyMode = $false

# The Czech sheet is no longer the active tab; its last selection becomes a
# full-sheet selection.
$czech.Cells.Select()

# Leave the new Swiss sheet selected with A16 as the active cell.
$newSheet.Range("A16").Select()
